$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 2745.8948
$ws.Range("I86").Value = 2312.6667
$ws.Range("J86").Value = 3488.5715
$ws.Range("K86").Value = 2312.6667
$ws.Range("L86").Value = 3488.5715
$ws.Range("M86").Value = -1189.6667
$ws.Range("N86").Value = -5734.5715

# Row 89
$ws.Range("H89").Value = 2745.8948
$ws.Range("I89").Value = 2312.6667
$ws.Range("J89").Value = 3488.5715
$ws.Range("K89").Value = 11563.3335
$ws.Range("L89").Value = 17442.8575
$ws.Range("M89").Value = -5947.333500000001
$ws.Range("N89").Value = -28674.8575

# Row 125
$ws.Range("H125").Value = 1313.6
$ws.Range("I125").Value = 1142
$ws.Range("K125").Value = 10278
$ws.Range("M125").Value = -7818

# Row 130
$ws.Range("H130").Value = 32500
$ws.Range("J130").Value = 32500
$ws.Range("L130").Value = 32500
$ws.Range("N130").Value = -42540

# Row 141
$ws.Range("H141").Value = 4997.9
$ws.Range("I141").Value = 2219.889
$ws.Range("J141").Value = 30000
$ws.Range("K141").Value = 6659.667
$ws.Range("L141").Value = 90000
$ws.Range("M141").Value = -1479.667
$ws.Range("N141").Value = -100360

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2760.5386
$ws.Range("I61").Value = 2760.5386
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2760.5386
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -2548.5386

# Row 97
$ws.Range("H97").Value = 823.7273
$ws.Range("J97").Value = 1168
$ws.Range("L97").Value = 1168
$ws.Range("N97").Value = -2160

# Row 102
$ws.Range("H102").Value = 1510.75
$ws.Range("I102").Value = 1370.2858
$ws.Range("J102").Value = 2494
$ws.Range("K102").Value = 1370.2858
$ws.Range("L102").Value = 2494
$ws.Range("M102").Value = 251.7141999999999
$ws.Range("N102").Value = -5738

# Row 104
$ws.Range("H104").Value = 41345
$ws.Range("J104").Value = 41345
$ws.Range("L104").Value = 41345
$ws.Range("N104").Value = -48333

# Row 135
$ws.Range("H135").Value = 69494.53999999999
$ws.Range("J135").Value = 69494.53999999999
$ws.Range("L135").Value = 69494.53999999999
$ws.Range("N135").Value = -79634.53999999999

# Row 136
$ws.Range("H136").Value = 2760.5386
$ws.Range("I136").Value = 2760.5386
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8281.6158
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -5731.6158

$ws = $wb.Worksheets.Item("BSM")
# Row 30
$ws.Range("H30").Value = 2500
$ws.Range("J30").Value = 3500
$ws.Range("L30").Value = 3500
$ws.Range("N30").Value = -3750

# Row 99
$ws.Range("H99").Value = 39989.184
$ws.Range("I99").Value = 39989.184
$ws.Range("K99").Value = 39989.184
$ws.Range("M99").Value = -38491.184

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1792.5454
$ws.Range("I31").Value = 1688.1177
$ws.Range("J31").Value = 2147.6
$ws.Range("K31").Value = 1688.1177
$ws.Range("L31").Value = 2147.6
$ws.Range("M31").Value = -1393.1177
$ws.Range("N31").Value = -2737.6

# Row 34
$ws.Range("H34").Value = 1792.5454
$ws.Range("I34").Value = 1688.1177
$ws.Range("J34").Value = 2147.6
$ws.Range("K34").Value = 1688.1177
$ws.Range("L34").Value = 2147.6
$ws.Range("M34").Value = -1486.1177
$ws.Range("N34").Value = -2551.6

# Row 62
$ws.Range("H62").Value = 6091.5
$ws.Range("I62").Value = 7385
$ws.Range("K62").Value = 7385
$ws.Range("M62").Value = -6761

# Row 65
$ws.Range("H65").Value = 6091.5
$ws.Range("I65").Value = 7385
$ws.Range("K65").Value = 36925
$ws.Range("M65").Value = -33805

# Row 94
$ws.Range("H94").Value = 9750.416999999999
$ws.Range("I94").Value = 17223.166
$ws.Range("K94").Value = 17223.166
$ws.Range("M94").Value = -16772.166

# Row 105
$ws.Range("H105").Value = 2495.9524
$ws.Range("I105").Value = 2613.1177
$ws.Range("K105").Value = 2613.1177
$ws.Range("M105").Value = -866.1176999999998

# Row 122
$ws.Range("H122").Value = 468091
$ws.Range("J122").Value = 6166.5835
$ws.Range("L122").Value = 18499.7505
$ws.Range("N122").Value = -23399.7505

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 777
$ws.Range("I122").Value = 1166.3334
$ws.Range("J122").Value = 582.3333
$ws.Range("K122").Value = 10497.0006
$ws.Range("L122").Value = 5240.9997
$ws.Range("M122").Value = -8047.000599999999
$ws.Range("N122").Value = -10140.9997

# Row 129
$ws.Range("H129").Value = 612.4
$ws.Range("I129").Value = 612.4
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 1837.2
$ws.Range("L129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = 3162.8

# Row 131
$ws.Range("H131").Value = 3134.7083
$ws.Range("I131").Value = 1959.3846
$ws.Range("J131").Value = 4523.727
$ws.Range("K131").Value = 5878.1538
$ws.Range("L131").Value = 13571.181
$ws.Range("M131").Value = -838.1538
$ws.Range("N131").Value = -23651.181

# Row 132
$ws.Range("H132").Value = 2997.5
$ws.Range("J132").Value = 2995
$ws.Range("L132").Value = 26955
$ws.Range("N132").Value = -32015

# Row 133
$ws.Range("H133").Value = 3990
$ws.Range("I133").Value = 3990
$ws.Range("K133").Value = 11970
$ws.Range("M133").Value = -6910

# Row 134
$ws.Range("H134").Value = 4606.636
$ws.Range("I134").Value = 3167.3
$ws.Range("K134").Value = 9501.900000000001
$ws.Range("M134").Value = -4431.900000000001

# Row 137
$ws.Range("H137").Value = 2659.182
$ws.Range("I137").Value = 2321.5715
$ws.Range("J137").Value = 3250
$ws.Range("K137").Value = 6964.7145
$ws.Range("L137").Value = 9750
$ws.Range("M137").Value = -1864.7145
$ws.Range("N137").Value = -19950

# Row 139
$ws.Range("H139").Value = 142858290
$ws.Range("I139").Value = 142858290
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 428574870
$ws.Range("L139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -428569730

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 1047.5714
$ws.Range("I132").Value = 655.5
$ws.Range("K132").Value = 1966.5
$ws.Range("M132").Value = 563.5

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1366
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 1488
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 1488
$ws.Range("M16").Value = -830
$ws.Range("N16").Value = -1828

# Row 93
$ws.Range("H93").Value = 8713.941000000001
$ws.Range("I93").Value = 8438.4
$ws.Range("J93").Value = 9107.571
$ws.Range("K93").Value = 8438.4
$ws.Range("L93").Value = 9107.571
$ws.Range("M93").Value = -7190.4
$ws.Range("N93").Value = -11603.571

# Row 105
$ws.Range("H105").Value = 36371.668
$ws.Range("J105").Value = 36371.668
$ws.Range("L105").Value = 36371.668
$ws.Range("N105").Value = -43359.668

$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 10711.857
$ws.Range("J45").Value = 11283
$ws.Range("L45").Value = 11283
$ws.Range("N45").Value = -12265

# Row 54
$ws.Range("H54").Value = 33133
$ws.Range("J54").Value = 30499.5
$ws.Range("L54").Value = 30499.5
$ws.Range("N54").Value = -31539.5

# Row 128
$ws.Range("H128").Value = 60000
$ws.Range("J128").Value = 60000
$ws.Range("L128").Value = 60000
$ws.Range("N128").Value = -69960

# Row 137
$ws.Range("H137").Value = 69791.586
$ws.Range("J137").Value = 69791.586
$ws.Range("L137").Value = 69791.586
$ws.Range("N137").Value = -79991.586
